# Apply the "Fixed update to excel issue" edit:
#  1. Rename "Requested quantity" header -> "Weekly_PO_Qty" on the "Weekly Quantity" sheet
#  2. Rename "Requested quantity" header -> "Monthly_PO_Qty" on the "Monthly Trend" sheet
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

# --- 1 & 2: rename the B1 headers on the existing sheets -------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" worksheet at the end ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$hdr = $wsForecast.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# Data rows (33 rows, r2:r34)
$data = New-Object 'object[,]' 33,4
$data[0,0] = 45025.99999999999; $data[0,1] = 48; $data[0,2] = -13.56049401143895; $data[0,3] = 104.4284552085926
$data[1,0] = 45032.99999999999; $data[1,1] = 48; $data[1,2] = -13.60625831117636; $data[1,3] = 110.3739844141056
$data[2,0] = 45039.99999999999; $data[2,1] = 48; $data[2,2] = -11.51174670484102; $data[2,3] = 110.6284745735361
$data[3,0] = 45046.99999999999; $data[3,1] = 49; $data[3,2] = -8.658053035974049; $data[3,3] = 108.563482427946
$data[4,0] = 45053.99999999999; $data[4,1] = 49; $data[4,2] = -10.2005209836438; $data[4,3] = 108.3971067566398
$data[5,0] = 45067.99999999999; $data[5,1] = 49; $data[5,2] = -8.057471487444175; $data[5,3] = 106.0747859240479
$data[6,0] = 45074.99999999999; $data[6,1] = 49; $data[6,2] = -7.865457783079977; $data[6,3] = 105.2186413413726
$data[7,0] = 45102.99999999999; $data[7,1] = 50; $data[7,2] = -12.58132578783702; $data[7,3] = 112.4275623819988
$data[8,0] = 45109.99999999999; $data[8,1] = 50; $data[8,2] = -11.50761553446732; $data[8,3] = 114.0969001016254
$data[9,0] = 45116.99999999999; $data[9,1] = 50; $data[9,2] = -8.251916747494327; $data[9,3] = 109.3559021432273
$data[10,0] = 45144.99999999999; $data[10,1] = 51; $data[10,2] = -11.6730580159299; $data[10,3] = 108.9085552183519
$data[11,0] = 45158.99999999999; $data[11,1] = 51; $data[11,2] = -12.55872926658225; $data[11,3] = 108.2944433486225
$data[12,0] = 45221.99999999999; $data[12,1] = 53; $data[12,2] = -6.646349817871691; $data[12,3] = 113.2198860543224
$data[13,0] = 45270.99999999999; $data[13,1] = 54; $data[13,2] = -6.762310870224626; $data[13,3] = 115.0842328843318
$data[14,0] = 45277.99999999999; $data[14,1] = 54; $data[14,2] = -8.28253523274585; $data[14,3] = 117.7032407624686
$data[15,0] = 45298.99999999999; $data[15,1] = 54; $data[15,2] = -9.166811361532696; $data[15,3] = 113.4065738506884
$data[16,0] = 45403.99999999999; $data[16,1] = 57; $data[16,2] = -0.7011652873776104; $data[16,3] = 116.077448457721
$data[17,0] = 45410.99999999999; $data[17,1] = 57; $data[17,2] = -0.6511797666898081; $data[17,3] = 121.675248317953
$data[18,0] = 45417.99999999999; $data[18,1] = 57; $data[18,2] = -6.375221090546185; $data[18,3] = 115.5139830323426
$data[19,0] = 45424.99999999999; $data[19,1] = 57; $data[19,2] = 1.265408696847387; $data[19,3] = 119.7152281271258
$data[20,0] = 45431.99999999999; $data[20,1] = 58; $data[20,2] = -4.073685946912916; $data[20,3] = 112.6510107314296
$data[21,0] = 45438.99999999999; $data[21,1] = 58; $data[21,2] = 0.2372937753690497; $data[21,3] = 122.6866022006237
$data[22,0] = 45445.99999999999; $data[22,1] = 58; $data[22,2] = -7.271251314073092; $data[22,3] = 119.226487941545
$data[23,0] = 45452.99999999999; $data[23,1] = 58; $data[23,2] = -1.805470392652321; $data[23,3] = 118.2570907067025
$data[24,0] = 45459.99999999999; $data[24,1] = 58; $data[24,2] = 0.7520191228424874; $data[24,3] = 119.0357386222691
$data[25,0] = 45466.99999999999; $data[25,1] = 58; $data[25,2] = -1.601054752831948; $data[25,3] = 119.4878099812896
$data[26,0] = 45473.99999999999; $data[26,1] = 59; $data[26,2] = -4.007280554549571; $data[26,3] = 116.7443818547925
$data[27,0] = 45480.99999999999; $data[27,1] = 59; $data[27,2] = -3.889204229201912; $data[27,3] = 116.7448164099437
$data[28,0] = 45487.99999999999; $data[28,1] = 59; $data[28,2] = 0.5327950195123483; $data[28,3] = 120.4246913846022
$data[29,0] = 45494.99999999999; $data[29,1] = 59; $data[29,2] = -3.609559763707881; $data[29,3] = 116.0035349726258
$data[30,0] = 45501.99999999999; $data[30,1] = 59; $data[30,2] = 2.294944991089726; $data[30,3] = 117.5329181919511
$data[31,0] = 45508.99999999999; $data[31,1] = 59; $data[31,2] = -4.118140962486004; $data[31,3] = 117.6286642315746
$data[32,0] = 45515.99999999999; $data[32,1] = 60; $data[32,2] = -1.471552288674176; $data[32,3] = 120.9673725575299

$rngData = $wsForecast.Range("A2:D34")
$rngData.Value = $data

# Column A data rows use the same date/time display format as the other sheets
$wsForecast.Range("A2:A34").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$wb.Save()
